$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.928.30"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.28"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.33"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4723"
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3676"
$ws.Range("E8").Value = "  +2.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07186"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9225"
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.60"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07603"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.830.27"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.308"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.395"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.41"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008640"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.947.33"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.55"
$ws.Range("E21").Value = "  +2.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.025"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.22"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.009"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.35"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.892"
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08855"
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.281"
$ws.Range("E31").Value = "  +4.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7480"
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.164"
$ws.Range("E33").Value = "  +3.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.781"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.491"
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05261"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01949"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.972"
$ws.Range("E39").Value = "  +2.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5206"
$ws.Range("E40").Value = "  +2.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.900"
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1512"
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.199"
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("E44").Value = "  +4.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4696"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.03"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.605"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.35"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8851"
